# "deleted video info in excel, added stymuls right to the tasks in a list format"
#
# Column B ("VideoIsNecessary") is no longer needed, so remove it entirely.
# Deleting the whole column shifts Text1/Text2/Text3 and the per-row stimulus
# lists (toy/candy/juice/etc. options) one column to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire "VideoIsNecessary" column, shifting everything to its
# right (Text1, Text2, Text3, and the stimulus lists) one column left.
$ws.Range("B:B").Delete()

# While tidying up the stimulus list, fix the "ябклуко" typo -> "яблуко".
# That cell used to be F2, and after the column shift above it is now E2.
$ws.Cells.Item(2, 5).Value = "яблуко"

# Leave the cursor where the author left it when saving.
$ws.Range("B10").Select()
